$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 137, shifting rows 137:219 down to 138:220
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new data record
$ws.Cells.Item(137, 1).Value = 4
$ws.Cells.Item(137, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(137, 3).Value = "Los Lagos"
$ws.Cells.Item(137, 4).Value = 44824
$ws.Cells.Item(137, 5).Value = 10
$ws.Cells.Item(137, 6).Value = 100112009
$ws.Cells.Item(137, 7).Value = "Acelga"
$ws.Cells.Item(137, 8).Value = "Sin especificar"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 80
$ws.Cells.Item(137, 11).Value = 1500
$ws.Cells.Item(137, 12).Value = 1500
$ws.Cells.Item(137, 13).Value = 1500
$ws.Cells.Item(137, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(137, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(137, 16).Value = 1000
$ws.Cells.Item(137, 17).Value = 1.5
$ws.Cells.Item(137, 18).Value = "Hortaliza"
